# Guiao_Apresentacao.docx - apply the textual edits described in the commit.
$d = $word.ActiveDocument

function FindRange($text) {
    $rng = $d.Content
    $null = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return $rng
}

function ReplaceText($oldText, $newText) {
    $null = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
}

# ---------------------------------------------------------------------------
# 1) "A tecnologia de casa inteligente, também conhecida como automação
#    residencial, fornece ..." -> wrap the word "automação " (soon to become
#    "automatização ") in the bookmark _Hlk512535789, same as upstream edit.
# ---------------------------------------------------------------------------
$rng = FindRange("automação ")
$d.Bookmarks.Add("_Hlk512535789", $rng)

# ---------------------------------------------------------------------------
# 2) Move the _GoBack bookmark: remove it from its old spot (in the
#    "...de modo que " / "sincronizar" paragraph) - it will be re-added
#    later, inside the "cafeteira" sentence, after that text is edited.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 3) All the straightforward word/phrase substitutions.
# ---------------------------------------------------------------------------

# "automação" -> "automatização" everywhere (9 occurrences total).
ReplaceText "automação" "automatização"

# " rápida olhada no smartphone" -> " vista de olhos rápida no smartphone"
ReplaceText "uma rápida olhada no smartphone" "uma vista de olhos rápida no smartphone"

# "... US $ 100 milhões construindo sua casa inteligente." -> "... na construção da sua casa inteligente."
ReplaceText "milhões construindo sua casa inteligente" "milhões na construção da sua casa inteligente"

# "uma miríade de outros dispositivos on-line. É tudo sobre a Internet das Coisas." ->
# "um grande número de outros dispositivos on-line. É tudo sobre a Internet das Coisas (IoT)."
ReplaceText "uma miríade de outros dispositivos on-line" "um grande número de outros dispositivos on-line"
ReplaceText "É tudo sobre a Internet das Coisas." "É tudo sobre a Internet das Coisas (IoT)."

# "A Internet das Coisas é uma expressão que se refere ... está ficando maior e
# melhor a cada dia. Todos os componentes eletrônicos da sua casa são um jogo
# justo para essa revolução tecnológica, da sua geladeira ao seu forno." ->
ReplaceText "uma expressão que se refere" "uma expressão, que se refere"
ReplaceText "está ficando maior e melhor a cada dia. Todos os componentes eletrônicos da sua casa são" "está a ficar cada dia maior e melhor. Todos os componentes eletrônicos de uma casa são"
ReplaceText "tecnológica, da sua geladeira ao seu forno." "tecnológica, do figorifico ao forno."

# "Algumas televisões inteligentes também incluem reconhecimento de voz ou gesto." ->
# "Algumas televisões inteligentes, também, incluem reconhecimento de voz ou gesto."
ReplaceText "televisões inteligentes também incluem reconhecimento de voz ou gesto." "televisões inteligentes, também, incluem reconhecimento de voz ou gesto."

# "As lâmpadas inteligentes também podem se regular" -> "... também se podem regular"
ReplaceText "também podem se regular com base na disponibilidade de luz do dia" "também se podem regular com base na disponibilidade de luz do dia"

# "lembrar os usuários de alterar os filtros" -> "lembrar os utilizadores de alterar os filtros"
ReplaceText "lembrar os usuários de alterar os filtros" "lembrar os utilizadores de alterar os filtros"

# "os moradores podem monitorar suas casas quando estão fora" -> "... monitorar as suas casas ..."
ReplaceText "podem monitorar suas casas quando estão fora" "podem monitorar as suas casas quando estão fora"

# "Plantas de casa e gramados podem ser regadas" -> "Plantas de casa e relvado podem ser regadas"
ReplaceText "Plantas de casa e gramados podem ser regadas" "Plantas de casa e relvado podem ser regadas"

# "alarme disparar; Geladeiras inteligentes" -> "alarme disparar; Figoríficos inteligentes"
ReplaceText "alarme disparar; Geladeiras inteligentes" "alarme disparar; Figoríficos inteligentes"

# "... desligar a água para que não haja uma inundação em seu porão." -> "... inundação."
ReplaceText "desligar a água para que não haja uma inundação em seu porão." "desligar a água para que não haja uma inundação."

# "... cafeteira esquecida deixada ou uma porta ..." -> "... cafeteira em funcionamento esquecida ou uma porta ..."
ReplaceText "cafeteira esquecida deixada ou uma porta" "cafeteira em funcionamento esquecida ou uma porta"

# "níveis de estresse" -> "níveis de stresse"
ReplaceText "níveis de estresse" "níveis de stresse"

# ---------------------------------------------------------------------------
# 4) Re-add the _GoBack bookmark as a collapsed range right after
#    "...cafeteira em funcionamento esquecida " and before "ou uma porta...".
# ---------------------------------------------------------------------------
$rng2 = FindRange("em funcionamento esquecida ")
$rng2.Collapse(0)  # wdCollapseEnd
$d.Bookmarks.Add("_GoBack", $rng2)
